# Apply updated crypto price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.993.23"
$ws.Range("E2").Value = "  -0.82%  "

$ws.Range("D3").Value = "3.387.43"

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'573.21"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").Value = "'137.19"
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.383.72"
$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").Value = "  -0.92%  "

$ws.Range("D10").Value = "'7.64"
$ws.Range("E10").Value = "  +1.85%  "

$ws.Range("E11").Value = "  -2.22%  "

$ws.Range("E12").Value = "  -2.59%  "

$ws.Range("D13").Value = "3.962.75"
$ws.Range("E13").Value = "  +0.20%  "

$ws.Range("E14").Value = "  -0.75%  "

$ws.Range("D15").Value = "'26.56"
$ws.Range("E15").Value = "  +2.31%  "

$ws.Range("D16").Value = "3.383.76"
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").Value = "'0.0000171"
$ws.Range("E17").Value = "  -2.97%  "

$ws.Range("D18").Value = "61.036.59"
$ws.Range("E18").Value = "  -0.86%  "

$ws.Range("D19").Value = "'13.87"
$ws.Range("E19").Value = "  -1.61%  "

$ws.Range("D20").Value = "'5.85"
$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("D21").Value = "'9.30"
$ws.Range("E21").Value = "  -0.98%  "

$ws.Range("D22").Value = "'374.95"
$ws.Range("E22").Value = "  -0.49%  "

$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "'0.551"
$ws.Range("E23").Value = "  -1.11%  "

$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.511.02"
$ws.Range("E24").Value = "  -0.29%  "

$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("D26").Value = "'70.81"
$ws.Range("E26").Value = "  -0.75%  "

$ws.Range("D27").Value = "'0.0000123"
$ws.Range("E27").Value = "  -2.74%  "

$ws.Range("E28").Value = "  -7.04%  "

$ws.Range("D29").Value = "'0.173"
$ws.Range("E29").Value = "  +7.38%  "

$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("D31").Value = "'7.36"

$ws.Range("D32").Value = "'8.07"
$ws.Range("E32").Value = "  -2.41%  "

$ws.Range("D33").Value = "'2.14"
$ws.Range("E33").Value = "  -1.74%  "

$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("D35").Value = "'23.38"
$ws.Range("E35").Value = "  -0.34%  "

$ws.Range("E36").Value = "  -3.22%  "

$ws.Range("D37").Value = "'1.55"
$ws.Range("E37").Value = "  +0.39%  "

$ws.Range("D38").Value = "'6.82"
$ws.Range("E38").Value = "  -0.17%  "

$ws.Range("D39").Value = "'164.76"
$ws.Range("E39").Value = "  -0.37%  "

$ws.Range("D40").Value = "'0.0766"
$ws.Range("E40").Value = "  -1.21%  "

$ws.Range("D41").Value = "'25.81"
$ws.Range("E41").Value = "  +4.01%  "

$ws.Range("E42").Value = "  +1.10%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("D45").Value = "'41.91"
$ws.Range("E45").Value = "  +0.99%  "

$ws.Range("D46").Value = "'4.37"
$ws.Range("E46").Value = "  -1.02%  "

$ws.Range("E47").Value = "  -4.17%  "

$ws.Range("D48").Value = "2.506.95"
$ws.Range("E48").Value = "  +6.87%  "

$ws.Range("D49").Value = "'23.61"
$ws.Range("E49").Value = "  +3.44%  "

$ws.Range("E50").Value = "  -1.04%  "

$ws.Range("E51").Value = "  +2.48%  "
